$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update Version value (row 3, column B) from 0.4.0 to 0.7.0
$ws.Range("B3").Value = "0.7.0"

# Remove the "Jurisdiction" / "Chile" row entirely (row 11), shifting all
# subsequent rows up by one.
$ws.Rows.Item(11).Delete()
